# Updates the cryptocurrency price/volume table to reflect the latest
# scrape performed by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.129.07'
$ws.Range("E2").Value = '  +1.77%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.066.35'
$ws.Range("E3").Value = '  +0.97%  '

# Row 4
$ws.Range("E4").Value = '  +0.41%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.47'
$ws.Range("E5").Value = '  +2.36%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.11'
$ws.Range("E6").Value = '  +1.88%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.30%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.064.58'
$ws.Range("E8").Value = '  +0.85%  '

# Row 9
$ws.Range("E9").Value = '  +4.82%  '

# Row 10
$ws.Range("E10").Value = '  +3.62%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.15'
$ws.Range("E11").Value = '  -10.82%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.493'
$ws.Range("E12").Value = '  +9.77%  '

# Row 13
$ws.Range("E13").Value = '  +4.11%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.74'
$ws.Range("E14").Value = '  +3.75%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.565.75'
$ws.Range("E15").Value = '  +2.01%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.141.55'
$ws.Range("E16").Value = '  +2.38%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.064.55'
$ws.Range("E17").Value = '  +1.83%  '

# Row 18
$ws.Range("E18").Value = '  +2.06%  '

# Row 19
$ws.Range("E19").Value = '  +2.36%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.89'
$ws.Range("E20").Value = '  +1.56%  '

# Row 21
$ws.Range("E21").Value = '  +4.04%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  +3.75%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.53'
$ws.Range("E23").Value = '  +14.01%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.60'
$ws.Range("E24").Value = '  +3.52%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.59'
$ws.Range("E25").Value = '  +3.43%  '

# Row 27
$ws.Range("E27").Value = '  +3.20%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.19'
$ws.Range("E28").Value = '  +5.41%  '

# Row 29
$ws.Range("E29").Value = '  +1.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.56%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.39'
$ws.Range("E31").Value = '  +2.63%  '

# Row 32
$ws.Range("E32").Value = '  +1.10%  '

# Row 33
$ws.Range("E33").Value = '  +3.31%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("E34").Value = '  +4.53%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.28'
$ws.Range("E35").Value = '  +6.00%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.93'
$ws.Range("E36").Value = '  -0.17%  '

# Row 37
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '453.63'
$ws.Range("E37").Value = '  -0.27%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0412'
$ws.Range("E38").Value = '  +4.56%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0820'
$ws.Range("E39").Value = '  +0.32%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.85'
$ws.Range("E40").Value = '  +9.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.036.16'
$ws.Range("E41").Value = '  +2.57%  '

# Row 42
$ws.Range("E42").Value = '  +1.60%  '

# Row 43
$ws.Range("E43").Value = '  +1.27%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.269'
$ws.Range("E44").Value = '  +7.22%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '28.04'
$ws.Range("E45").Value = '  +3.26%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.24'
$ws.Range("E46").Value = '  +10.95%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.113'
$ws.Range("E48").Value = '  +3.11%  '

# Row 49
$ws.Range("E49").Value = '  +1.95%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '118.47'
$ws.Range("E50").Value = '  +2.94%  '

# Row 51
$ws.Range("E51").Value = '  +3.44%  '
